# Update the "Case locations and outbreaks / public exposure sites" table
# with the latest published exposure-site list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set: header row plus the 12 new/updated exposure-site
# rows reflecting the latest DHHS update (the previous 4 data rows are gone).
$header = @("Location", "Site", "Exposure period", "Notes", "Exist")
$data = @(
    @("Abbotsford",       "Bodriggy Brewing Company  245 Johnston Street, Abbotsford VIC 3067", "28/12/20 2:50pm-5:30pm",     "Case dined at venue",          "new"),
    @("Brighton",         "Brighton Beach",                                                      "26/12/20 12pm - 1pm",        "Case attended beach",          "old"),
    @("Docklands",        "Melbourne Boat Hire - Yarra River Cruise. 45 Newquay Promenade, Docklands VIC 3008", "28/12/2020 11.26am-2:00pm", "Case attended venue", "new"),
    @("Hampton",          "Merrymen Cafe, 2 Small Street, Hampton VIC",                          "28-12-2020 1:00pm-2:00pm",   "Case ate in store",            "new"),
    @("Hampton",          "Merrymen Cafe, 2 Small Street, Hampton VIC",                          "28-12-2020 1:20pm-2:30pm",   "Case ate in store",            "old"),
    @("Hampton",          "Merrymen Cafe, 2 Small Street, Hampton VIC",                          "28-12-2020 1:30pm-2:30pm",   "Case ate in store",            "new"),
    @("McKinnon",         "260 McKinnon Road, McKinnon VIC 3204",                                "23-12-2020 4:00pm-6:00pm",   "Case had hair cut in store",   "old"),
    @("Melbourne",        "Melbourne Central Lion Hotel, 211 La Trobe Street",                   "28/12/2020 10pm - 12.30am",  "Case attended venue",          "new"),
    @("Moorabin",         "Grape and Grain Liquor Cellars, 14/16 Station St",                    "21/12/20 2pm - 10pm  22/12/20 10am - 6pm  24/12/20 1pm - 10pm  28/12/20 8.05pm - 8.47pm  29/12/20 12pm - 4pm", "Cases workplace", "new"),
    @("Sandringham Line", "Metro Train line Sandringham",                                        "28/12/20 7pm -7.50pm",       "Travelled by train from Sandringham Station to Parliament Station", "new"),
    @("Southbank",        "Left Bank Melbourne Restaurant and Cocktail Bar  1 Southbank Boulevard, Southbank", "25/12/20 12:00pm-02:30pm", "Case attended bar", "old"),
    @("Southbank",        "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank",  "23/12/20 8:00pm-10:00pm",    "Case attended restaurant",     "old")
)

# Clear the whole used range (including the header) before re-writing, so
# the shared-string table is rebuilt fresh in top-to-bottom / left-to-right
# scan order instead of preserving stale indices for reused strings.
$ws.Cells.ClearContents()

for ($c = 0; $c -lt $header.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $header[$c]
}

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Resize columns: A, B, D, E auto-fit to their (new, longer) content; C
# ("Exposure period") keeps a generous manually-set width to fit the
# multi-line Moorabin exposure list.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 96.5

# Leave the selection where the author last left it (last site's "Site" cell).
$ws.Range("B12").Select()
